# The workbook "trains_and_teams" sheet has 300 data rows (A1:D300).
# Column D currently holds "None" for every row (shared string 305).
# The target edit replaces column D with a repeating 22-value cycle:
#   AT1, AT2, ..., AT20, ST1, ST2, AT1, AT2, ...
# which introduces 22 brand-new shared strings (appended right after
# "None" in sharedStrings.xml, in first-seen order), taking uniqueCount
# from 306 to 328.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    "AT1","AT2","AT3","AT4","AT5","AT6","AT7","AT8","AT9","AT10",
    "AT11","AT12","AT13","AT14","AT15","AT16","AT17","AT18","AT19","AT20",
    "ST1","ST2"
)

for ($r = 1; $r -le 300; $r++) {
    $idx = ($r - 1) % $values.Length
    $ws.Cells.Item($r, 4).Value = $values[$idx]
}

# Match the updated view state: the sheet was scrolled down a bit further
# and the selection moved from H298 to D301:D308.
$win = $excel.ActiveWindow
$ws.Range("D301:D308").Select()
$win.ScrollRow = 283
$win.ScrollColumn = 1

Write-Host "Updated column D (D1:D300) with AT/ST cycle and refreshed selection"
